$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 24
$ws.Range("E6").Value = 46
$ws.Range("E7").Value = 23
$ws.Range("F7").Value = 12
$ws.Range("H7").Value = 12
$ws.Range("E10").Value = 19
$ws.Range("F10").Value = 7
$ws.Range("H10").Value = 7
$ws.Range("F16").Value = 77
$ws.Range("H16").Value = 77
$ws.Range("E18").Value = 79
